# Auto-generated edit script: updates cryptos list values/percentages
# per commit "Updated cryptos list on Mon Jun  3 15:26:34 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "69.201.86"
$ws.Range("E2").Value = "  +1.67%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.786.54"
$ws.Range("E3").Value = "  -0.21%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5: BNB
$ws.Range("D5").Value = "'629.69"
$ws.Range("E5").Value = "  +4.69%  "

# Row 6: Solana
$ws.Range("D6").Value = "'164.31"
$ws.Range("E6").Value = "  -0.42%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "3.784.85"
$ws.Range("E7").Value = "  -0.23%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.08%  "

# Row 9: XRP
$ws.Range("E9").Value = "  +0.44%  "

# Row 10: Dogecoin
$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  +1.15%  "

# Row 11: Cardano
$ws.Range("E11").Value = "  +0.07%  "

# Row 12: Toncoin
$ws.Range("D12").Value = "'6.63"
$ws.Range("E12").Value = "  +2.32%  "

# Row 13: ShibaInu
$ws.Range("E13").Value = "  -0.72%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "'35.42"
$ws.Range("E14").Value = "  -1.09%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.420.17"
$ws.Range("E15").Value = "  -0.30%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "3.795.76"
$ws.Range("E16").Value = "  +0.20%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "69.257.38"
$ws.Range("E17").Value = "  +1.72%  "

# Row 18: Chainlink
$ws.Range("D18").Value = "'17.97"
$ws.Range("E18").Value = "  -2.24%  "

# Row 19: Polkadot
$ws.Range("E19").Value = "  -0.01%  "

# Row 20: TRON
$ws.Range("E20").Value = "  -1.20%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'468.44"
$ws.Range("E21").Value = "  +1.51%  "

# Row 22: Uniswap
$ws.Range("E22").Value = "  -0.62%  "

# Row 23: Polygon
$ws.Range("D23").Value = "'0.705"
$ws.Range("E23").Value = "  +0.58%  "

# Row 24: PEPE
$ws.Range("D24").Value = "'0.0000150"
$ws.Range("E24").Value = "  +1.12%  "

# Row 25: Litecoin
$ws.Range("D25").Value = "'83.23"
$ws.Range("E25").Value = "  +0.22%  "

# Row 26: InternetComputer(DFINITY)
$ws.Range("D26").Value = "'12.08"
$ws.Range("E26").Value = "  +0.39%  "

# Row 27: Fetch.AI
$ws.Range("E27").Value = "  +1.73%  "

# Row 28: RenderToken
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.02"
$ws.Range("E28").Value = "  +0.04%  "

# Row 29: Dai
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30: WrappedeETH
$ws.Range("D30").Value = "3.929.99"
$ws.Range("E30").Value = "  -0.40%  "

# Row 31: PancakeSwap
$ws.Range("D31").Value = "'2.68"
$ws.Range("E31").Value = "  +0.99%  "

# Row 32: ImmutableX
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.23"
$ws.Range("E32").Value = "  -0.15%  "

# Row 33: NEARProtocol
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'7.30"
$ws.Range("E33").Value = "  -0.69%  "

# Row 34: EthereumClassic
$ws.Range("D34").Value = "'29.01"
$ws.Range("E34").Value = "  -1.12%  "

# Row 35: Binance-PegBSC-USD
$ws.Range("E35").Value = "  +0.02%  "

# Row 36: Aptos
$ws.Range("D36").Value = "'9.02"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37: RenzoRestakedETH
$ws.Range("D37").Value = "3.733.50"
$ws.Range("E37").Value = "  -0.29%  "

# Row 38: Hedera
$ws.Range("E38").Value = "  +2.54%  "

# Row 39: Kaspa
$ws.Range("D39").Value = "'0.150"
$ws.Range("E39").Value = "  +8.10%  "

# Row 40: dogwifhat
$ws.Range("D40").Value = "'3.32"
$ws.Range("E40").Value = "  +0.30%  "

# Row 41: Filecoin
$ws.Range("D41").Value = "'5.84"
$ws.Range("E41").Value = "  -0.20%  "

# Row 42: Mantle
$ws.Range("E42").Value = "  -1.90%  "

# Row 43: FirstDigitalUSD
$ws.Range("E43").Value = "  -0.08%  "

# Row 44: USDe
$ws.Range("E44").Value = "  +0.04%  "

# Row 45: TheGraph
$ws.Range("D45").Value = "'0.299"
$ws.Range("E45").Value = "  -0.18%  "

# Row 46: Monero
$ws.Range("D46").Value = "'153.10"
$ws.Range("E46").Value = "  +1.11%  "

# Row 47: Stacks
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.92"
$ws.Range("E47").Value = "  +3.00%  "

# Row 48: OKB
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'46.86"
$ws.Range("E48").Value = "  -1.25%  "

# Row 49: Arweave
$ws.Range("D49").Value = "'42.68"
$ws.Range("E49").Value = "  -1.53%  "

# Row 50: Cosmos
$ws.Range("D50").Value = "'8.44"
$ws.Range("E50").Value = "  +0.92%  "

# Row 51: ONDO
$ws.Range("D51").Value = "'1.39"
$ws.Range("E51").Value = "  +2.74%  "

